$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-12-21 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-12-22 Sunday", 2) | Out-Null

# Update the 20x5 table of arithmetic expressions
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "14+11="
$tbl.Cell(1, 2).Range.Text = "50-32="
$tbl.Cell(1, 3).Range.Text = "96-77="
$tbl.Cell(1, 4).Range.Text = "24+74="
$tbl.Cell(1, 5).Range.Text = "72-44="

$tbl.Cell(2, 1).Range.Text = "54-30="
$tbl.Cell(2, 2).Range.Text = "74-25="
$tbl.Cell(2, 3).Range.Text = "25+18="
$tbl.Cell(2, 4).Range.Text = "91-11="
$tbl.Cell(2, 5).Range.Text = "95-51="

$tbl.Cell(3, 1).Range.Text = "8+38="
$tbl.Cell(3, 2).Range.Text = "50+49="
$tbl.Cell(3, 3).Range.Text = "45+22="
$tbl.Cell(3, 4).Range.Text = "54+40="
$tbl.Cell(3, 5).Range.Text = "8+24="

$tbl.Cell(4, 1).Range.Text = "22+16="
$tbl.Cell(4, 2).Range.Text = "47+43="
$tbl.Cell(4, 3).Range.Text = "2+17="
$tbl.Cell(4, 4).Range.Text = "87-33="
$tbl.Cell(4, 5).Range.Text = "30+0="

$tbl.Cell(5, 1).Range.Text = "58+2="
$tbl.Cell(5, 2).Range.Text = "37+43="
$tbl.Cell(5, 3).Range.Text = "33-19="
$tbl.Cell(5, 4).Range.Text = "27+56="
$tbl.Cell(5, 5).Range.Text = "67-49="

$tbl.Cell(6, 1).Range.Text = "63-1="
$tbl.Cell(6, 2).Range.Text = "68+14="
$tbl.Cell(6, 3).Range.Text = "33+39="
$tbl.Cell(6, 4).Range.Text = "66-21="
$tbl.Cell(6, 5).Range.Text = "87-69="

$tbl.Cell(7, 1).Range.Text = "47+44="
$tbl.Cell(7, 2).Range.Text = "77-28="
$tbl.Cell(7, 3).Range.Text = "17+66="
$tbl.Cell(7, 4).Range.Text = "39+5="
$tbl.Cell(7, 5).Range.Text = "12+35="

$tbl.Cell(8, 1).Range.Text = "70-34="
$tbl.Cell(8, 2).Range.Text = "86-6="
$tbl.Cell(8, 3).Range.Text = "25-11="
$tbl.Cell(8, 4).Range.Text = "34+18="
$tbl.Cell(8, 5).Range.Text = "93-90="

$tbl.Cell(9, 1).Range.Text = "43-14="
$tbl.Cell(9, 2).Range.Text = "63-15="
$tbl.Cell(9, 3).Range.Text = "24+40="
$tbl.Cell(9, 4).Range.Text = "46-3="
$tbl.Cell(9, 5).Range.Text = "82-61="

$tbl.Cell(10, 1).Range.Text = "57-32="
$tbl.Cell(10, 2).Range.Text = "80+1="
$tbl.Cell(10, 3).Range.Text = "72-7="
$tbl.Cell(10, 4).Range.Text = "57-1="
$tbl.Cell(10, 5).Range.Text = "66+24="

$tbl.Cell(11, 1).Range.Text = "2+16="
$tbl.Cell(11, 2).Range.Text = "23-8="
$tbl.Cell(11, 3).Range.Text = "51+19="
$tbl.Cell(11, 4).Range.Text = "37+20="
$tbl.Cell(11, 5).Range.Text = "81-48="

$tbl.Cell(12, 1).Range.Text = "99-19="
$tbl.Cell(12, 2).Range.Text = "38+18="
$tbl.Cell(12, 3).Range.Text = "70+9="
$tbl.Cell(12, 4).Range.Text = "32+6="
$tbl.Cell(12, 5).Range.Text = "78-62="

$tbl.Cell(13, 1).Range.Text = "48+13="
$tbl.Cell(13, 2).Range.Text = "20+26="
$tbl.Cell(13, 3).Range.Text = "68+19="
$tbl.Cell(13, 4).Range.Text = "83+15="
$tbl.Cell(13, 5).Range.Text = "41+44="

$tbl.Cell(14, 1).Range.Text = "43-27="
$tbl.Cell(14, 2).Range.Text = "28-22="
$tbl.Cell(14, 3).Range.Text = "46+5="
$tbl.Cell(14, 4).Range.Text = "40+53="
$tbl.Cell(14, 5).Range.Text = "77-51="

$tbl.Cell(15, 1).Range.Text = "64-53="
$tbl.Cell(15, 2).Range.Text = "11+81="
$tbl.Cell(15, 3).Range.Text = "49-24="
$tbl.Cell(15, 4).Range.Text = "16+52="
$tbl.Cell(15, 5).Range.Text = "20+75="

$tbl.Cell(16, 1).Range.Text = "65-10="
$tbl.Cell(16, 2).Range.Text = "50+15="
$tbl.Cell(16, 3).Range.Text = "43+47="
$tbl.Cell(16, 4).Range.Text = "35-3="
$tbl.Cell(16, 5).Range.Text = "40-31="

$tbl.Cell(17, 1).Range.Text = "74-18="
$tbl.Cell(17, 2).Range.Text = "13+20="
$tbl.Cell(17, 3).Range.Text = "41+5="
$tbl.Cell(17, 4).Range.Text = "37+25="
$tbl.Cell(17, 5).Range.Text = "27+9="

$tbl.Cell(18, 1).Range.Text = "32+14="
$tbl.Cell(18, 2).Range.Text = "90-58="
$tbl.Cell(18, 3).Range.Text = "63+31="
$tbl.Cell(18, 4).Range.Text = "77-64="
$tbl.Cell(18, 5).Range.Text = "70-17="

$tbl.Cell(19, 1).Range.Text = "41+52="
$tbl.Cell(19, 2).Range.Text = "92+6="
$tbl.Cell(19, 3).Range.Text = "94-43="
$tbl.Cell(19, 4).Range.Text = "26+19="
$tbl.Cell(19, 5).Range.Text = "83-58="

$tbl.Cell(20, 1).Range.Text = "39+50="
$tbl.Cell(20, 2).Range.Text = "16+6="
$tbl.Cell(20, 3).Range.Text = "82-28="
$tbl.Cell(20, 4).Range.Text = "62+31="
$tbl.Cell(20, 5).Range.Text = "2+80="
